# Auto-generated Excel COM-interop script
# Applies scheduled market-price refresh values to the Leve profit tables
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 933.6818
$ws.Range("I33").Value = 507.4737
$ws.Range("J33").Value = 3633
$ws.Range("K33").Value = 507.4737
$ws.Range("L33").Value = 3633
$ws.Range("M33").Value = -278.4737
$ws.Range("N33").Value = -4091
$ws.Range("H58").Value = 1735.6666
$ws.Range("I58").Value = 254.375
$ws.Range("J58").Value = 3428.5715
$ws.Range("K58").Value = 763.125
$ws.Range("L58").Value = 10285.7145
$ws.Range("M58").Value = -613.125
$ws.Range("N58").Value = -10585.7145
$ws.Range("H125").Value = 2494.6667
$ws.Range("I125").Value = 932
$ws.Range("J125").Value = 2807.2
$ws.Range("K125").Value = 8388
$ws.Range("L125").Value = 25264.8
$ws.Range("M125").Value = -5928
$ws.Range("N125").Value = -30184.8
$ws.Range("H132").Value = 4030.1667
$ws.Range("I132").Value = 3919.6177
$ws.Range("K132").Value = 11758.8531
$ws.Range("M132").Value = -9228.8531

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H61").Value = 1745.1765
$ws.Range("I61").Value = 1275.9584
$ws.Range("J61").Value = 2871.3
$ws.Range("K61").Value = 1275.9584
$ws.Range("L61").Value = 2871.3
$ws.Range("M61").Value = -1063.9584
$ws.Range("N61").Value = -3295.3
$ws.Range("H64").Value = 24658.334
$ws.Range("J64").Value = 24658.334
$ws.Range("L64").Value = 24658.334
$ws.Range("N64").Value = -25154.334
$ws.Range("H67").Value = 24658.334
$ws.Range("J67").Value = 24658.334
$ws.Range("L67").Value = 24658.334
$ws.Range("N67").Value = -26374.334
$ws.Range("H136").Value = 1745.1765
$ws.Range("I136").Value = 1275.9584
$ws.Range("J136").Value = 2871.3
$ws.Range("K136").Value = 3827.8752
$ws.Range("L136").Value = 8613.900000000001
$ws.Range("M136").Value = -1277.8752
$ws.Range("N136").Value = -13713.9

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 31491.666
$ws.Range("J62").Value = 31491.666
$ws.Range("L62").Value = 31491.666
$ws.Range("N62").Value = -32863.666
$ws.Range("H65").Value = 31491.666
$ws.Range("J65").Value = 31491.666
$ws.Range("L65").Value = 94474.998
$ws.Range("N65").Value = -101338.998
$ws.Range("H105").Value = 2438.25
$ws.Range("I105").Value = 2663.3333
$ws.Range("J105").Value = 2386.3076
$ws.Range("K105").Value = 2663.3333
$ws.Range("L105").Value = 2386.3076
$ws.Range("M105").Value = -916.3332999999998
$ws.Range("N105").Value = -5880.3076

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 10371.429
$ws.Range("J80").Value = 10371.429
$ws.Range("L80").Value = 10371.429
$ws.Range("N80").Value = -12617.429
$ws.Range("H83").Value = 10371.429
$ws.Range("J83").Value = 10371.429
$ws.Range("L83").Value = 31114.287
$ws.Range("N83").Value = -42346.287
$ws.Range("H87").Value = 42450
$ws.Range("J87").Value = 42450
$ws.Range("L87").Value = 42450
$ws.Range("N87").Value = -44822
$ws.Range("H90").Value = 42450
$ws.Range("J90").Value = 42450
$ws.Range("L90").Value = 127350
$ws.Range("N90").Value = -139206
$ws.Range("H94").Value = 166667940
$ws.Range("I94").Value = 1000000000
$ws.Range("K94").Value = 1000000000
$ws.Range("M94").Value = -999999549
$ws.Range("H132").Value = 2109.476
$ws.Range("I132").Value = 1194.2941
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 3582.8823
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -1052.8823
$ws.Range("N132").Value = -23057
$ws.Range("H134").Value = 2719.3333
$ws.Range("I134").Value = 2541.3547
$ws.Range("J134").Value = 3822.8
$ws.Range("K134").Value = 7624.0641
$ws.Range("L134").Value = 11468.4
$ws.Range("M134").Value = -5089.0641
$ws.Range("N134").Value = -16538.4

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 10000
$ws.Range("J46").Value = 10000
$ws.Range("L46").Value = 30000
$ws.Range("N46").Value = -30182
$ws.Range("H107").Value = 1027.421
$ws.Range("I107").Value = 258.75
$ws.Range("J107").Value = 1232.4
$ws.Range("K107").Value = 776.25
$ws.Range("L107").Value = 3697.2
$ws.Range("M107").Value = 1143.75
$ws.Range("N107").Value = -7537.200000000001

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2318.6572
$ws.Range("I126").Value = 2519.7368
$ws.Range("J126").Value = 2079.875
$ws.Range("K126").Value = 7559.2104
$ws.Range("L126").Value = 6239.625
$ws.Range("M126").Value = -5089.2104
$ws.Range("N126").Value = -11179.625

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3817.75
$ws.Range("I7").Value = 3769.0908
$ws.Range("J7").Value = 3877.2222
$ws.Range("K7").Value = 3769.0908
$ws.Range("L7").Value = 3877.2222
$ws.Range("M7").Value = -3657.0908
$ws.Range("N7").Value = -4101.2222
$ws.Range("H16").Value = 1317.2667
$ws.Range("I16").Value = 1366.2307
$ws.Range("J16").Value = 999
$ws.Range("K16").Value = 1366.2307
$ws.Range("L16").Value = 999
$ws.Range("M16").Value = -1196.2307
$ws.Range("N16").Value = -1339
$ws.Range("H40").Value = 3368.5151
$ws.Range("I40").Value = 3584.2778
$ws.Range("J40").Value = 3109.6
$ws.Range("K40").Value = 3584.2778
$ws.Range("L40").Value = 3584.2778
$ws.Range("M40").Value = -3448.2778
$ws.Range("N40").Value = -3381.6
$ws.Range("H122").Value = 3608.611
$ws.Range("I122").Value = 2626.6667
$ws.Range("K122").Value = 7880.000100000001
$ws.Range("M122").Value = -5430.000100000001
$ws.Range("H123").Value = 31500
$ws.Range("I123").Value = 30000
$ws.Range("K123").Value = 30000
$ws.Range("M123").Value = -25100
$ws.Range("H126").Value = 3817.75
$ws.Range("I126").Value = 3769.0908
$ws.Range("J126").Value = 3877.2222
$ws.Range("K126").Value = 11307.2724
$ws.Range("L126").Value = 11631.6666
$ws.Range("M126").Value = -8837.2724
$ws.Range("N126").Value = -16571.6666
$ws.Range("H132").Value = 1495261.9
$ws.Range("I132").Value = 2987546
$ws.Range("J132").Value = 2977.762
$ws.Range("K132").Value = 8962638
$ws.Range("L132").Value = 8933.286
$ws.Range("M132").Value = -8960108
$ws.Range("N132").Value = -13993.286
$ws.Range("H136").Value = 5440105
$ws.Range("I136").Value = 8340741.5
$ws.Range("J136").Value = 1412.5
$ws.Range("K136").Value = 25022224.5
$ws.Range("L136").Value = 4237.5
$ws.Range("M136").Value = -25019674.5
$ws.Range("N136").Value = -9337.5

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 71429910
$ws.Range("I122").Value = 142858020
$ws.Range("J122").Value = 1814.2858
$ws.Range("K122").Value = 428574060
$ws.Range("L122").Value = 5442.857400000001
$ws.Range("M122").Value = -428571610
$ws.Range("N122").Value = -10342.8574
$ws.Range("H126").Value = 3340.75
$ws.Range("I126").Value = 5041.2856
$ws.Range("J126").Value = 960
$ws.Range("K126").Value = 15123.8568
$ws.Range("L126").Value = 2880
$ws.Range("M126").Value = -12653.8568
$ws.Range("N126").Value = -7820
$ws.Range("H132").Value = 1655.2858
$ws.Range("I132").Value = 1361.6786
$ws.Range("J132").Value = 2829.7144
$ws.Range("K132").Value = 4085.0358
$ws.Range("L132").Value = 8489.143199999999
$ws.Range("M132").Value = -1555.0358
$ws.Range("N132").Value = -13549.1432
